# "De buttons aan StartScene toegevoegd"
#
# Logs a new entry (row 11) on the "week 50" sheet: a fifth StartScene
# button-related activity, with its begin/end time filled in. The
# existing shared formula in column G, and the downstream SUM/reference
# formulas (week 50!G18 and Totaal!B9), recalculate automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("week 50")

# Begin time (11:30) and end time (12:05) for the new activity row.
$ws.Range("C11").Value = 0.47916666666666669
$ws.Range("D11").Value = 0.50347222222222221

# Description of the activity -> lands in sharedStrings as a new entry.
$ws.Range("F11").Value = "Alle buttons aangemaakt en gedrawed op het StartScene."

# Row grows to fit the wrapped activity text, same as the other
# multi-line rows (8 and 10) above it.
$ws.Rows.Item(11).RowHeight = 28.5

# Selection moves up one row after the edit.
$ws.Range("F12").Select() | Out-Null
